# Apply crypto price/volume updates + two row label swaps (rows 36/37 and 39/40)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values below are stored as text in the source workbook (e.g. "1.00", "311.09").
# Excel auto-converts such strings to numbers on plain assignment, so force text format,
# assign the value, then restore the default (Normal) style so no stray formatting remains.

$ws.Range('D2').Value = '41.308.79'
$ws.Range('E2').Value = '  -2.99%  '
$ws.Range('D3').Value = '2.460.43'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.77%  '
$ws.Range('E7').Value = '  -3.32%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.68%  '
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.06%  '
$ws.Range('D14').Value = '2.839.93'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').Value = '2.445.67'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.784'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').Value = '41.229.68'
$ws.Range('E18').Value = '  -3.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.71%  '
$ws.Range('D20').Value = '0.0₃0916'
$ws.Range('E20').Value = '  -3.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.65%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.85%  '
$ws.Range('E28').Value = '  -4.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.68'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.53%  '
$ws.Range('E33').Value = '  -6.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.59'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0753'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.36%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.84%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.89'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.19%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.114'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.87%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.102'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.82%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '1.958.80'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0284'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '69.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.75'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.179'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.80%  '
